$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new blank "test" sheet between "Merchandise" and "Paids".
#    It gets the header row only (same header as the item sheets).
# ---------------------------------------------------------------------------
$paids = $wb.Worksheets.Item("Paids")
$testSheet = $wb.Worksheets.Add($paids)
$testSheet.Name = "test"

$testSheet.Range("A1").Value = "Items: "
$testSheet.Range("B1").Value = "Regular"
$testSheet.Range("D1").Value = "Modifiers:"
$testSheet.Range("E1").Value = "Confirmed Completed: True"

# ---------------------------------------------------------------------------
# 2. "Bottled Beer": Troegs IPA price 3.5 -> 3.1, "Truly" renamed to "test"
#    (price cleared, moved down to row 6), "Twisted Tea" row removed.
# ---------------------------------------------------------------------------
$bottledBeer = $wb.Worksheets.Item("Bottled Beer")
$bottledBeer.Range("B3").Value = 3.1
$bottledBeer.Range("A4").Value = "test"
$bottledBeer.Range("B4").ClearContents()
$bottledBeer.Rows(5).Delete()
$bottledBeer.Range("A6").Value = $bottledBeer.Range("A4").Value2
$bottledBeer.Range("A4").ClearContents()

# ---------------------------------------------------------------------------
# 3. "Mixed Drinks": remove the "Trash Can" row entirely.
# ---------------------------------------------------------------------------
$mixedDrinks = $wb.Worksheets.Item("Mixed Drinks")
$mixedDrinks.Rows(10).Delete()

# ---------------------------------------------------------------------------
# 4. "Employees": every employee's PIN becomes "1234" (kept as text, like
#    the original zero-padded PIN codes, rather than being auto-converted
#    to a number).
# ---------------------------------------------------------------------------
$employees = $wb.Worksheets.Item("Employees")
$employees.Range("B2:B12").NumberFormat = "@"
for ($r = 2; $r -le 12; $r++) {
    $employees.Cells.Item($r, 2).Value = "1234"
}
